$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: L2 is empty, but carries the same style as K2 (thick-bottom border row) ---
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

# --- Row 3: L3 = 2022, same style as K3 (year header) ---
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L3").Value = 2022

# --- Row 4: L4 = 370, same style as K4 ---
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 370

# --- Row 5: L5 = 137, same style as K5 ---
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 137

# --- Row 6: L6 = 314, same style as K6 ---
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = 314

# --- Row 7: L7 = 121, same style as K7 ---
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = 121

# --- Row 8: L8 = 50, same base style as K8 but with an explicit "#,##0" number format
#     (this produces the new cellXfs entry seen in the diff) ---
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").NumberFormat = "#,##0"
$ws.Range("L8").Value = 50

# --- Row 9: L9 = 16, same style as K9 ---
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L9").Value = 16

# --- Update the selected cell to L2, matching the saved view state ---
$null = $ws.Range("L2").Select()
